$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing transformer/RNN parameter values in row 2
$ws.Range("C2").Value = 34
$ws.Range("D2").Value = 110
$ws.Range("F2").Value = 0.0001
$ws.Range("H2").Value = 10

# Add new transformer parameter columns (headers + values)
$ws.Range("K1").Value = "d_model"
$ws.Range("L1").Value = "num_layers"
$ws.Range("K2").Value = 16
$ws.Range("L2").Value = 1

# Reflect the active cell/selection left by the author
$ws.Range("J10").Select()
